$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add column L, mirroring column K (2020 / 6.18), including formatting.
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)

$ws.Range("L3").Value = 2020
$ws.Range("L4").Value = 6.18

# Leave the selection on M12, matching the saved view state.
$ws.Range("M12").Select()
